$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.837.92'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '3.216.08'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('E4').Value = '  -0.02%  '
$r = $ws.Range('D5')
$r.NumberFormat = "@"
$r.Value = '592.79'
$r.ClearFormats()
$ws.Range('E5').Value = '  -1.16%  '
$r = $ws.Range('D6')
$r.NumberFormat = "@"
$r.Value = '137.71'
$r.ClearFormats()
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.214.91'
$ws.Range('E8').Value = '  -1.39%  '
$r = $ws.Range('D9')
$r.NumberFormat = "@"
$r.Value = '0.515'
$r.ClearFormats()
$ws.Range('E9').Value = '  +1.00%  '
$r = $ws.Range('D10')
$r.NumberFormat = "@"
$r.Value = '0.143'
$r.ClearFormats()
$ws.Range('E10').Value = '  -2.19%  '
$ws.Range('E11').Value = '  -2.11%  '
$r = $ws.Range('D12')
$r.NumberFormat = "@"
$r.Value = '0.457'
$r.ClearFormats()
$ws.Range('E12').Value = '  -0.57%  '
$r = $ws.Range('D13')
$r.NumberFormat = "@"
$r.Value = '0.0000241'
$r.ClearFormats()
$ws.Range('E13').Value = '  -0.45%  '
$r = $ws.Range('D14')
$r.NumberFormat = "@"
$r.Value = '35.41'
$r.ClearFormats()
$ws.Range('E14').Value = '  +4.36%  '
$ws.Range('D15').Value = '3.745.45'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('D17').Value = '3.216.27'
$ws.Range('E17').Value = '  -1.41%  '
$ws.Range('D18').Value = '63.847.95'
$ws.Range('E18').Value = '  +1.02%  '
$r = $ws.Range('D19')
$r.NumberFormat = "@"
$r.Value = '6.59'
$r.ClearFormats()
$ws.Range('E19').Value = '  -1.91%  '
$r = $ws.Range('D20')
$r.NumberFormat = "@"
$r.Value = '465.53'
$r.ClearFormats()
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('E21').Value = '  +1.96%  '
$ws.Range('E22').Value = '  -2.29%  '
$r = $ws.Range('D23')
$r.NumberFormat = "@"
$r.Value = '7.70'
$r.ClearFormats()
$ws.Range('E23').Value = '  -1.63%  '
$r = $ws.Range('D24')
$r.NumberFormat = "@"
$r.Value = '13.45'
$r.ClearFormats()
$ws.Range('E24').Value = '  -0.67%  '
$r = $ws.Range('D25')
$r.NumberFormat = "@"
$r.Value = '83.57'
$r.ClearFormats()
$ws.Range('E25').Value = '  -0.42%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  -1.07%  '
$ws.Range('E28').Value = '  -0.04%  '
$r = $ws.Range('D29')
$r.NumberFormat = "@"
$r.Value = '7.88'
$r.ClearFormats()
$ws.Range('E29').Value = '  -0.81%  '
$r = $ws.Range('D31')
$r.NumberFormat = "@"
$r.Value = '6.87'
$r.ClearFormats()
$ws.Range('E31').Value = '  -2.64%  '
$r = $ws.Range('D32')
$r.NumberFormat = "@"
$r.Value = '27.59'
$r.ClearFormats()
$ws.Range('E32').Value = '  -1.34%  '
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('E35').Value = '  -3.21%  '
$r = $ws.Range('D36')
$r.NumberFormat = "@"
$r.Value = '5.92'
$r.ClearFormats()
$ws.Range('E36').Value = '  +0.43%  '
$r = $ws.Range('D37')
$r.NumberFormat = "@"
$r.Value = '51.69'
$r.ClearFormats()
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('D38').Value = '0.0₃0735'
$ws.Range('E38').Value = '  +1.83%  '
$r = $ws.Range('D39')
$r.NumberFormat = "@"
$r.Value = '0.0395'
$r.ClearFormats()
$ws.Range('E39').Value = '  +0.38%  '
$r = $ws.Range('D40')
$r.NumberFormat = "@"
$r.Value = '2.77'
$r.ClearFormats()
$ws.Range('E40').Value = '  +4.60%  '
$r = $ws.Range('D41')
$r.NumberFormat = "@"
$r.Value = '404.68'
$r.ClearFormats()
$ws.Range('E41').Value = '  -4.42%  '
$ws.Range('E42').Value = '  -0.55%  '
$r = $ws.Range('D43')
$r.NumberFormat = "@"
$r.Value = '0.113'
$r.ClearFormats()
$ws.Range('E43').Value = '  -2.65%  '
$ws.Range('D44').Value = '2.837.85'
$ws.Range('E44').Value = '  -7.58%  '
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('E46').Value = '  +0.36%  '
$r = $ws.Range('D47')
$r.NumberFormat = "@"
$r.Value = '129.20'
$r.ClearFormats()
$ws.Range('E47').Value = '  +2.16%  '
$r = $ws.Range('D48')
$r.NumberFormat = "@"
$r.Value = '36.05'
$r.ClearFormats()
$ws.Range('E48').Value = '  +0.66%  '
$r = $ws.Range('D50')
$r.NumberFormat = "@"
$r.Value = '25.86'
$r.ClearFormats()
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('E51').Value = '  -0.20%  '
